{"js": "// Replace each two-digit multiplication expression in the document with\n// its updated counterpart, in document order. Each source string is\n// unique within the document, so a simple search + replace per pair is\n// unambiguous.\nconst replacements = [\n  [\"21\u00d720=\", \"78\u00d759=\"],\n  [\"37\u00d719=\", \"57\u00d712=\"],\n  [\"82\u00d789=\", \"87\u00d734=\"],\n  [\"47\u00d736=\", \"52\u00d759=\"],\n  [\"72\u00d714=\", \"16\u00d733=\"],\n  [\"86\u00d754=\", \"54\u00d749=\"],\n  [\"14\u00d792=\", \"96\u00d754=\"],\n  [\"68\u00d779=\", \"28\u00d757=\"],\n  [\"22\u00d770=\", \"41\u00d794=\"],\n  [\"62\u00d772=\", \"83\u00d743=\"],\n  [\"26\u00d775=\", \"25\u00d743=\"],\n  [\"37\u00d788=\", \"17\u00d795=\"],\n  [\"57\u00d715=\", \"24\u00d735=\"],\n  [\"51\u00d776=\", \"96\u00d798=\"],\n  [\"12\u00d722=\", \"67\u00d716=\"],\n  [\"98\u00d753=\", \"93\u00d773=\"],\n  [\"14\u00d769=\", \"99\u00d775=\"],\n  [\"55\u00d749=\", \"83\u00d726=\"],\n  [\"74\u00d740=\", \"39\u00d763=\"],\n  [\"40\u00d768=\", \"51\u00d712=\"],\n  [\"29\u00d791=\", \"73\u00d728=\"],\n  [\"51\u00d789=\", \"60\u00d737=\"],\n  [\"58\u00d738=\", \"68\u00d739=\"],\n  [\"89\u00d750=\", \"81\u00d799=\"],\n  [\"55\u00d748=\", \"79\u00d711=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first match \u2014 each source string occurs exactly\n  // once, so this also guards against accidentally touching a\n  // previously-inserted replacement value.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression in the document with\n# its updated counterpart. Each source string occurs exactly once in the\n# document, so Find/Replace with MatchCase and a single replacement per\n# pair is unambiguous and preserves the existing run formatting.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"21\u00d720=\"; New = \"78\u00d759=\" },\n    @{ Old = \"37\u00d719=\"; New = \"57\u00d712=\" },\n    @{ Old = \"82\u00d789=\"; New = \"87\u00d734=\" },\n    @{ Old = \"47\u00d736=\"; New = \"52\u00d759=\" },\n    @{ Old = \"72\u00d714=\"; New = \"16\u00d733=\" },\n    @{ Old = \"86\u00d754=\"; New = \"54\u00d749=\" },\n    @{ Old = \"14\u00d792=\"; New = \"96\u00d754=\" },\n    @{ Old = \"68\u00d779=\"; New = \"28\u00d757=\" },\n    @{ Old = \"22\u00d770=\"; New = \"41\u00d794=\" },\n    @{ Old = \"62\u00d772=\"; New = \"83\u00d743=\" },\n    @{ Old = \"26\u00d775=\"; New = \"25\u00d743=\" },\n    @{ Old = \"37\u00d788=\"; New = \"17\u00d795=\" },\n    @{ Old = \"57\u00d715=\"; New = \"24\u00d735=\" },\n    @{ Old = \"51\u00d776=\"; New = \"96\u00d798=\" },\n    @{ Old = \"12\u00d722=\"; New = \"67\u00d716=\" },\n    @{ Old = \"98\u00d753=\"; New = \"93\u00d773=\" },\n    @{ Old = \"14\u00d769=\"; New = \"99\u00d775=\" },\n    @{ Old = \"55\u00d749=\"; New = \"83\u00d726=\" },\n    @{ Old = \"74\u00d740=\"; New = \"39\u00d763=\" },\n    @{ Old = \"40\u00d768=\"; New = \"51\u00d712=\" },\n    @{ Old = \"29\u00d791=\"; New = \"73\u00d728=\" },\n    @{ Old = \"51\u00d789=\"; New = \"60\u00d737=\" },\n    @{ Old = \"58\u00d738=\"; New = \"68\u00d739=\" },\n    @{ Old = \"89\u00d750=\"; New = \"81\u00d799=\" },\n    @{ Old = \"55\u00d748=\"; New = \"79\u00d711=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $found = $range.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($pair.Old)\"\n    }\n}\n"}
